# Season-record columns: Wins / Losses / Ties appended after the existing
# team-stat columns (A:AC), replicated across every player row (2-39).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Header row (row 1): new headers in AD1:AF1, matching the style of the
# existing header cells (copy format from AC1, the last existing header).
$headers = @("Wins", "Losses", "Ties")
$headerCols = @(30, 31, 32)  # AD, AE, AF

for ($i = 0; $i -lt $headers.Length; $i++) {
    $col = $headerCols[$i]
    $cell = $ws.Cells.Item(1, $col)
    $cell.Value = $headers[$i]
}

# Copy the header style from an existing header cell (AC1) onto the new
# header cells so they match (bold, centered, bordered).
$ws.Range("AC1").Copy() | Out-Null
$ws.Range("AD1:AF1").PasteSpecial(-4122) | Out-Null

# Re-set the values in case PasteSpecial(formats) touched anything.
$ws.Cells.Item(1, 30).Value = "Wins"
$ws.Cells.Item(1, 31).Value = "Losses"
$ws.Cells.Item(1, 32).Value = "Ties"

# --- Data rows (2-39): season record is identical for every player on the
# roster (team-wide W-L-T), so fill straight down.
for ($r = 2; $r -le 39; $r++) {
    $ws.Cells.Item($r, 30).Value = 96
    $ws.Cells.Item($r, 31).Value = 66
    $ws.Cells.Item($r, 32).Value = 0
}
